{"js": "// The document opens with the Title, Author and Abstract paragraphs each\n// split word-by-word across many runs (one run per word, one run per\n// separating space). The edit simply merges each of those paragraphs\n// down to a single run holding the full paragraph text, leaving the\n// paragraph style / every other paragraph untouched.\n\nfunction runOoxmlForParagraph(styleId, text) {\n  // Build a minimal, well-formed OOXML \"flat OPC\" package snippet that\n  // describes one paragraph with one run so we can hand it to\n  // `insertOoxml(..., \"Replace\")`. Using OOXML (rather than\n  // `insertText`) lets us control the exact `<w:t>` serialization\n  // (keeping `xml:space=\"preserve\"`) the same way the source document\n  // already writes every text run.\n  const escaped = text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    '<w:pPr><w:pStyle w:val=\"' +\n    styleId +\n    '\"/></w:pPr>' +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">\" +\n    escaped +\n    \"</w:t></w:r>\" +\n    \"</w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst targets = [\n  { style: \"Title\", text: \"Answers: Solving exponential equations\" },\n  {\n    style: \"Author\",\n    text: \"Zo\u00eb Gemmell, Isabella Lewis, Akshat Srivastava\",\n  },\n  {\n    style: \"Abstract\",\n    text: \"Answers to questions relating to solving exponential equations.\",\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const item of paragraphs.items) {\n  item.load(\"style\");\n}\nawait context.sync();\n\nfor (const target of targets) {\n  const paragraph = paragraphs.items.find((p) => p.style === target.style);\n  if (!paragraph) {\n    continue;\n  }\n  paragraph.insertOoxml(\n    runOoxmlForParagraph(target.style, target.text),\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# The document opens with the Title, Author and Abstract paragraphs each\n# split word-by-word across many runs (one run per word, one run per\n# separating space). The edit simply merges each of those paragraphs\n# down to a single run holding the full paragraph text, leaving the\n# paragraph style / every other paragraph untouched.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText {\n    param(\n        $Paragraph,\n        $NewText\n    )\n\n    $range = $Paragraph.Range\n    # Range.Text includes the trailing paragraph mark (chr 13); trim it so\n    # Find only ever targets the paragraph's visible text.\n    $oldText = $range.Text.TrimEnd([char]13)\n\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $NewText\n\n    # wdReplaceAll = 2, wdFindContinue = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nSet-ParagraphText $d.Paragraphs(1) \"Answers: Solving exponential equations\"\nSet-ParagraphText $d.Paragraphs(2) \"Zo\u00eb Gemmell, Isabella Lewis, Akshat Srivastava\"\nSet-ParagraphText $d.Paragraphs(4) \"Answers to questions relating to solving exponential equations.\"\n"}
